# Update the "想去人数" (want-to-go count) figures in column F for the
# rows that changed between scrapes, on both sheets that carry the full
# data table: "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1699
    $ws.Range("F6").Value = 472
    $ws.Range("F9").Value = 628
    $ws.Range("F10").Value = 410
}
